$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None

MSG: The committee did not reach a decision regarding which movie to show on Friday.
"
$ws.Range("D2").Value = "no_decision, "
$ws.Range("C3").Value = "MSG: None

MSG: The decision regarding Friday's movie was not reached, so I have called the no_decision function.
"
$ws.Range("D3").Value = "no_decision, "
$ws.Range("C4").Value = "MSG: None

MSG: The decision has been recorded as no selection being made for the movie on Friday.
"
$ws.Range("D4").Value = "no_decision, "
$ws.Range("C5").Value = "MSG: None

MSG: The decision has been recorded. The movie ""Barbie"" will be shown on Friday.
"
$ws.Range("D5").Value = "Barbie_was_selected, "
$ws.Range("C6").Value = "MSG: None

MSG: The decision regarding Friday's movie has been recorded as ""no decision.""
"
$ws.Range("D6").Value = "no_decision, "
$ws.Range("C7").Value = "MSG: None

MSG: The decision has been recorded: there was no agreement on a movie choice for Friday.
"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None

MSG: The decision process concluded without selecting a movie for Friday, and thus no further action is taken.
"
$ws.Range("D8").Value = "no_decision, "
$ws.Range("C9").Value = "MSG: None

MSG: The decision has been recorded, indicating that no definitive choice about Friday's movie has been made.
"
$ws.Range("D9").Value = "no_decision, "
$ws.Range("C10").Value = "MSG: None

MSG: The function has been called successfully, indicating that no decision was made regarding the movie to be shown on Friday.
"
$ws.Range("D10").Value = "no_decision, "
$ws.Range("C11").Value = "MSG: None

MSG: The decision has been recorded as no movie being selected for Friday.
"
$ws.Range("D11").Value = "no_decision, "
$ws.Range("C12").Value = "MSG: None

MSG: The decision about which movie to show on Friday cannot be made.
"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None

MSG: The decision has been recorded as having no outcome for the movie selection.
"
$ws.Range("D13").Value = "no_decision, "
$ws.Range("C14").Value = "MSG: None

MSG: The decision to acquire the rights for ""Barbie"" has been recorded.
"
$ws.Range("D14").Value = "Barbie_was_selected, "
$ws.Range("C15").Value = "MSG: None

MSG: The decision regarding the movie for Friday has been recorded as no decision being made.
"
$ws.Range("D15").Value = "no_decision, "
$ws.Range("C16").Value = "MSG: None

MSG: The decision has been recorded as no movie being selected for Friday.
"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None

MSG: No movie was selected in this meeting.
"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None

MSG: The decision has been recorded, indicating that no consensus was reached on a movie to show on Friday.
"
$ws.Range("D18").Value = "no_decision, "
$ws.Range("C19").Value = "MSG: None

MSG: The decision process did not result in an agreement on a movie to be shown on Friday; hence, no further action will be taken.
"
$ws.Range("D19").Value = "no_decision, "
$ws.Range("C20").Value = "MSG: None

MSG: The decision has been recorded as no decision regarding Friday's movie was made.
"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None

MSG: The decision-making process did not yield a final choice regarding the movie to be shown on Friday.
"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None

MSG: The decision has been recorded as no decision regarding Friday's movie was made.
"
$ws.Range("D22").Value = "no_decision, "
$ws.Range("C23").Value = "MSG: None

MSG: The decision has been made that no movie will be shown on Friday.
"
$ws.Range("D23").Value = "no_decision, "
$ws.Range("C24").Value = "MSG: None

MSG: The decision has been recorded as no movie was selected in this meeting.
"
$ws.Range("D24").Value = "no_decision, "
$ws.Range("C25").Value = "MSG: None

MSG: The decision has been recorded as no decision was made regarding Friday's movie.
"
$ws.Range("D25").Value = "no_decision, "
$ws.Range("C26").Value = "MSG: None

MSG: The decision to acquire the rights for ""Barbie"" has been successfully recorded.
"
$ws.Range("D26").Value = "Barbie_was_selected, "
$ws.Range("C27").Value = "MSG: None

MSG: The rights for both movies have been successfully acquired.
"
$ws.Range("D27").Value = "both_movies, "
$ws.Range("C28").Value = "MSG: None

MSG: The decision regarding which movie to show on Friday was not finalized, resulting in no movie being selected.
"
$ws.Range("D28").Value = "no_decision, "
$ws.Range("C29").Value = "MSG: None

MSG: The decision process concluded without reaching a consensus on which movie to show on Friday, so no specific movie rights will be acquired at this time.
"
$ws.Range("D29").Value = "no_decision, "
$ws.Range("C30").Value = "MSG: None

MSG: The decision to acquire the rights for both movies has been recorded successfully.
"
$ws.Range("D30").Value = "both_movies, "
$ws.Range("C31").Value = "MSG: None

MSG: The movie ""Barbie"" has been successfully selected for the assembly on Friday.
"
$ws.Range("D31").Value = "Barbie_was_selected, "
$ws.Range("C32").Value = "MSG: None

MSG: The rights to both movies have been successfully acquired for the screening.
"
$ws.Range("D32").Value = "both_movies, "
$ws.Range("C33").Value = "MSG: None

MSG: I have called the no_decision function, indicating that no decision was reached regarding the movie to be shown on Friday.
"
$ws.Range("D33").Value = "no_decision, "
$ws.Range("C34").Value = "MSG: None

MSG: The decision process did not lead to an agreement on which movie to show on Friday, so no action on acquiring rights is necessary.
"
$ws.Range("D34").Value = "no_decision, "
$ws.Range("C35").Value = "MSG: None

MSG: I have recorded the decision to acquire the rights for ""Barbie"" as the movie to be shown on Friday.
"
$ws.Range("D35").Value = "Barbie_was_selected, "
$ws.Range("C36").Value = "MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday was made.
"
$ws.Range("D36").Value = "no_decision, "
$ws.Range("C37").Value = "MSG: None

MSG: The decision has been recorded with no movie selected for Friday.
"
$ws.Range("D37").Value = "no_decision, "
$ws.Range("C38").Value = "MSG: None

MSG: The decision for Friday's movie cannot be made, as the committee ended the conversation without a definitive choice.
"
$ws.Range("D38").Value = "no_decision, "
$ws.Range("C39").Value = "MSG: None

MSG: The decision has been recorded successfully for ""Barbie"" to be shown on Friday.
"
$ws.Range("D39").Value = "Barbie_was_selected, "
$ws.Range("C40").Value = "MSG: None

MSG: The decision regarding Friday's movie has resulted in no conclusion.
"
$ws.Range("D40").Value = "no_decision, "
$ws.Range("C41").Value = "MSG: None

MSG: The decision has been recorded, and there is currently no selected movie for the Friday showing.
"
$ws.Range("D41").Value = "no_decision, "
$ws.Range("C42").Value = "MSG: None

MSG: The committee did not reach a decision regarding the movie to show on Friday. Therefore, the outcome is recorded as no decision made.
"
$ws.Range("D42").Value = "no_decision, "
$ws.Range("C43").Value = "MSG: None

MSG: The decision has been successfully recorded to acquire the rights for ""Barbie.""
"
$ws.Range("D43").Value = "Barbie_was_selected, "
$ws.Range("C44").Value = "MSG: None

MSG: The decision about what movie to show on Friday has not been reached.
"
$ws.Range("D44").Value = "no_decision, "
$ws.Range("C45").Value = "MSG: None

MSG: The rights to both movies have been acquired.
"
$ws.Range("D45").Value = "both_movies, "
$ws.Range("C46").Value = "MSG: None

MSG: A decision about which movie to show on Friday could not be reached.
"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None

MSG: The decision about Friday's movie was not reached.
"
$ws.Range("D47").Value = "no_decision, "
$ws.Range("C48").Value = "MSG: None

MSG: The decision to acquire rights for ""Oppenheimer"" has been recorded successfully.
"
$ws.Range("D48").Value = "Oppenheimer_was_selected, "
$ws.Range("C49").Value = "MSG: None

MSG: The decision has been recorded to acquire the rights for the movie ""Barbie.""
"
$ws.Range("D49").Value = "Barbie_was_selected, "
$ws.Range("C50").Value = "MSG: None

MSG: The decision has been made to acquire the rights for ""Barbie.""
"
$ws.Range("D50").Value = "Barbie_was_selected, "
$ws.Range("C51").Value = "MSG: None

MSG: The decision has been made to acquire the rights for both movies, ""Oppenheimer"" and ""Barbie.""
"
$ws.Range("D51").Value = "both_movies, "
$ws.Range("C52").Value = "MSG: None

MSG: The decision has been recorded, and rights for both movies will be acquired.
"
$ws.Range("D52").Value = "both_movies, "
$ws.Range("C53").Value = "MSG: None

MSG: The decision has been recorded as no_decision.
"
$ws.Range("D53").Value = "no_decision, "
$ws.Range("C54").Value = "MSG: None

MSG: The decision has been recorded, indicating that no movie was selected for Friday.
"
$ws.Range("D54").Value = "no_decision, "
$ws.Range("C55").Value = "MSG: None

MSG: The decision resulted in no agreement on a movie for Friday, and thus no movie rights will be acquired.
"
$ws.Range("D55").Value = "no_decision, "
$ws.Range("C56").Value = "MSG: None

MSG: The decision has been recorded, and no movie was selected for Friday.
"
$ws.Range("D56").Value = "no_decision, "
$ws.Range("C57").Value = "MSG: None

MSG: The decision has been made to acquire the rights for ""Oppenheimer.""
"
$ws.Range("D57").Value = "Oppenheimer_was_selected, "
$ws.Range("C58").Value = "MSG: None

MSG: I have executed the function to indicate that no decision was made regarding Friday's movie.
"
$ws.Range("D58").Value = "no_decision, "
$ws.Range("C59").Value = "MSG: None

MSG: The decision has been made to acquire the rights for both movies.
"
$ws.Range("D59").Value = "both_movies, "
$ws.Range("C60").Value = "MSG: None

MSG: The decision regarding the movie to be shown on Friday remains undecided.
"
$ws.Range("D60").Value = "no_decision, "
$ws.Range("C61").Value = "MSG: None

MSG: The decision indicates that no definitive choice was reached regarding which movie to show on Friday.
"
$ws.Range("D61").Value = "no_decision, "
$ws.Range("C62").Value = "MSG: None

MSG: The rights to both movies have been successfully acquired for showing on Friday.
"
$ws.Range("D62").Value = "both_movies, "
$ws.Range("C63").Value = "MSG: None

MSG: The decision regarding the movie to be shown on Friday resulted in no consensus, and therefore, no movie was selected.
"
$ws.Range("D63").Value = "no_decision, "
$ws.Range("C64").Value = "MSG: None

MSG: The decision about which movie to show on Friday could not be made, as the conversation ended without a clear agreement.
"
$ws.Range("D64").Value = "no_decision, "
$ws.Range("C65").Value = "MSG: None

MSG: The decision has been recorded, and no movie will be shown on Friday as there was no agreement reached by the committee.
"
$ws.Range("D65").Value = "no_decision, "
$ws.Range("C66").Value = "MSG: None

MSG: The decision to acquire the rights for ""Barbie"" has been recorded.
"
$ws.Range("D66").Value = "Barbie_was_selected, "
$ws.Range("C67").Value = "MSG: None

MSG: I have decided to acquire the rights to ""Barbie"" as the movie to be shown on Friday.
"
$ws.Range("D67").Value = "Barbie_was_selected, "

Write-Output "Updated rows 2-67 classifications"